$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.103.68"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").Value = "2.222.23"
$ws.Range("E3").Value = "  -0.49%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "291.80"
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "87.84"
$ws.Range("E6").Value = "  +2.00%  "
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +1.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "30.36"
$ws.Range("E10").Value = "  -0.57%  "
$ws.Range("E11").Value = "  -2.08%  "
$ws.Range("E12").Value = "  +3.46%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.48"
$ws.Range("E13").Value = "  +1.53%  "
$ws.Range("D14").Value = "2.569.60"
$ws.Range("E14").Value = "  -0.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "13.95"
$ws.Range("E15").Value = "  -1.34%  "
$ws.Range("D16").Value = "2.224.91"
$ws.Range("E16").Value = "  -1.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.727"
$ws.Range("E17").Value = "  +0.47%  "
$ws.Range("D18").Value = "40.062.95"
$ws.Range("E18").Value = "  +0.55%  "
$ws.Range("E20").Value = "  +7.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.82"
$ws.Range("E21").Value = "  +0.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.62"
$ws.Range("E22").Value = "  +0.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.62"
$ws.Range("E23").Value = "  +0.62%  "
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.47"
$ws.Range("E25").Value = "  +1.82%  "
$ws.Range("E26").Value = "  -0.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.68"
$ws.Range("E27").Value = "  -0.95%  "
$ws.Range("E28").Value = "  -4.43%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.21"
$ws.Range("E29").Value = "  -0.52%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "156.61"
$ws.Range("E30").Value = "  +1.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "31.75"
$ws.Range("E31").Value = "  -6.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").Value = "  -0.07%  "
$ws.Range("E33").Value = "  +2.05%  "
$ws.Range("E34").Value = "  +1.15%  "
$ws.Range("E35").Value = "  -1.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.88"
$ws.Range("E36").Value = "  +6.79%  "
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "15.65"
$ws.Range("E38").Value = "  -5.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0979"
$ws.Range("E39").Value = "  -2.15%  "
$ws.Range("E40").Value = "  +1.90%  "
$ws.Range("D41").Value = "2.116.29"
$ws.Range("E41").Value = "  +8.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.86"
$ws.Range("E42").Value = "  +2.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.14"
$ws.Range("E43").Value = "  -2.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "18.03"
$ws.Range("E44").Value = "  +11.27%  "
$ws.Range("E45").Value = "  -0.93%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.91"
$ws.Range("E46").Value = "  +3.41%  "
$ws.Range("E47").Value = "  +2.98%  "
$ws.Range("D48").Value = "2.436.25"
$ws.Range("E48").Value = "  -0.51%  "
$ws.Range("E49").Value = "  -0.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "69.51"
$ws.Range("E50").Value = "  -1.83%  "
$ws.Range("E51").Value = "  +2.90%  "